$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-33 ---
$ws.Range("B2").Value = "NSE:AARVEEDEN"
$ws.Range("C2").Value = "NSE:ANGELONE"
$ws.Range("D2").Value = "NSE:DIXON"
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()

$ws.Range("B3").Value = "NSE:ABSLNN50ET"
$ws.Range("C3").Value = "NSE:HINDMOTORS"
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()

$ws.Range("B4").Value = "NSE:AKI"
$ws.Range("C4").Value = "NSE:MARUTI"
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()

$ws.Range("B5").Value = "NSE:ASTRAZEN"
$ws.Range("C5").Value = "NSE:MITTAL"
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()

$ws.Range("B6").Value = "NSE:AXISHCETF"
$ws.Range("C6").Value = "NSE:ORIENTALTL"
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()

$ws.Range("B7").Value = "NSE:BEL"
$ws.Range("C7").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("F7").ClearContents()

$ws.Range("B8").Value = "NSE:BIKAJI"
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()

$ws.Range("B9").Value = "NSE:COCHINSHIP"
$ws.Range("C9").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("F9").ClearContents()

$ws.Range("B10").Value = "NSE:DCMSRIND"
$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("F10").ClearContents()

$ws.Range("B11").Value = "NSE:DCW"
$ws.Range("C11").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("F11").ClearContents()

$ws.Range("B12").Value = "NSE:DCXINDIA"
$ws.Range("C12").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("F12").ClearContents()

$ws.Range("B13").Value = "NSE:DELHIVERY"
$ws.Range("C13").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("F13").ClearContents()

$ws.Range("B14").Value = "NSE:DHANI"
$ws.Range("C14").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("E14").ClearContents()
$ws.Range("F14").ClearContents()

$ws.Range("B15").Value = "NSE:DTIL"
$ws.Range("C15").ClearContents()
$ws.Range("D15").ClearContents()
$ws.Range("E15").ClearContents()
$ws.Range("F15").ClearContents()

$ws.Range("B16").Value = "NSE:EMAMIPAP"
$ws.Range("C16").ClearContents()
$ws.Range("D16").ClearContents()
$ws.Range("E16").ClearContents()
$ws.Range("F16").ClearContents()

$ws.Range("B17").Value = "NSE:ESAFSFB"
$ws.Range("C17").ClearContents()
$ws.Range("D17").ClearContents()
$ws.Range("E17").ClearContents()
$ws.Range("F17").ClearContents()

$ws.Range("B18").Value = "NSE:GENUSPAPER"
$ws.Range("C18").ClearContents()
$ws.Range("D18").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("F18").ClearContents()

$ws.Range("B19").Value = "NSE:GRSE"
$ws.Range("C19").ClearContents()
$ws.Range("D19").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Range("F19").ClearContents()

$ws.Range("B20").Value = "NSE:HAL"
$ws.Range("C20").ClearContents()
$ws.Range("D20").ClearContents()
$ws.Range("E20").ClearContents()
$ws.Range("F20").ClearContents()

$ws.Range("B21").Value = "NSE:INOXGREEN"
$ws.Range("C21").ClearContents()
$ws.Range("D21").ClearContents()
$ws.Range("E21").ClearContents()
$ws.Range("F21").ClearContents()

$ws.Range("B22").Value = "NSE:IVP"
$ws.Range("C22").ClearContents()
$ws.Range("D22").ClearContents()
$ws.Range("E22").ClearContents()
$ws.Range("F22").ClearContents()

$ws.Range("B23").Value = "NSE:JSWINFRA"
$ws.Range("C23").ClearContents()
$ws.Range("D23").ClearContents()
$ws.Range("E23").ClearContents()
$ws.Range("F23").ClearContents()

$ws.Range("B24").Value = "NSE:KEC"
$ws.Range("C24").ClearContents()
$ws.Range("D24").ClearContents()
$ws.Range("E24").ClearContents()
$ws.Range("F24").ClearContents()

$ws.Range("B25").Value = "NSE:KINGFA"
$ws.Range("C25").ClearContents()
$ws.Range("D25").ClearContents()
$ws.Range("E25").ClearContents()
$ws.Range("F25").ClearContents()

$ws.Range("B26").Value = "NSE:KMSUGAR"
$ws.Range("C26").ClearContents()
$ws.Range("D26").ClearContents()
$ws.Range("E26").ClearContents()
$ws.Range("F26").ClearContents()

$ws.Range("B27").Value = "NSE:KOHINOOR"
$ws.Range("C27").ClearContents()
$ws.Range("D27").ClearContents()
$ws.Range("E27").ClearContents()
$ws.Range("F27").ClearContents()

$ws.Range("B28").Value = "NSE:KRBL"
$ws.Range("C28").ClearContents()
$ws.Range("D28").ClearContents()
$ws.Range("E28").ClearContents()
$ws.Range("F28").ClearContents()

$ws.Range("B29").Value = "NSE:KSHITIJPOL"
$ws.Range("C29").ClearContents()
$ws.Range("D29").ClearContents()
$ws.Range("E29").ClearContents()
$ws.Range("F29").ClearContents()

$ws.Range("B30").Value = "NSE:KUANTUM"
$ws.Range("C30").ClearContents()
$ws.Range("D30").ClearContents()
$ws.Range("E30").ClearContents()
$ws.Range("F30").ClearContents()

$ws.Range("B31").Value = "NSE:MANINDS"
$ws.Range("C31").ClearContents()
$ws.Range("D31").ClearContents()
$ws.Range("E31").ClearContents()
$ws.Range("F31").ClearContents()

$ws.Range("B32").Value = "NSE:MCLEODRUSS"
$ws.Range("C32").ClearContents()
$ws.Range("D32").ClearContents()
$ws.Range("E32").ClearContents()
$ws.Range("F32").ClearContents()

$ws.Range("B33").Value = "NSE:MIDHANI"
$ws.Range("C33").ClearContents()
$ws.Range("D33").ClearContents()
$ws.Range("E33").ClearContents()
$ws.Range("F33").ClearContents()

# --- Add new rows 34-42 (copy style of row 33 column A for the index style) ---
$ws.Range("A33").Copy($ws.Range("A34"))
$ws.Range("A34").Value = 32
$ws.Range("B34").Value = "NSE:MTARTECH"

$ws.Range("A33").Copy($ws.Range("A35"))
$ws.Range("A35").Value = 33
$ws.Range("B35").Value = "NSE:NAHARPOLY"

$ws.Range("A33").Copy($ws.Range("A36"))
$ws.Range("A36").Value = 34
$ws.Range("B36").Value = "NSE:NIRAJ"

$ws.Range("A33").Copy($ws.Range("A37"))
$ws.Range("A37").Value = 35
$ws.Range("B37").Value = "NSE:NUVOCO"

$ws.Range("A33").Copy($ws.Range("A38"))
$ws.Range("A38").Value = 36
$ws.Range("B38").Value = "NSE:OAL"

$ws.Range("A33").Copy($ws.Range("A39"))
$ws.Range("A39").Value = 37
$ws.Range("B39").Value = "NSE:ONEPOINT"

$ws.Range("A33").Copy($ws.Range("A40"))
$ws.Range("A40").Value = 38
$ws.Range("B40").Value = "NSE:PARAS"

$ws.Range("A33").Copy($ws.Range("A41"))
$ws.Range("A41").Value = 39
$ws.Range("B41").Value = "NSE:PTCIL"

$ws.Range("A33").Copy($ws.Range("A42"))
$ws.Range("A42").Value = 40
$ws.Range("B42").Value = "NSE:RML"
